$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, add row 354 by copying the constant (categorical) columns from row 353,
# since the new row inherits the same Mercado/Categoria/etc. as the rest of the block.
$ws.Cells.Item(354, 1).Value = $ws.Cells.Item(353, 1).Value2
$ws.Cells.Item(354, 2).Value = $ws.Cells.Item(353, 2).Value2
$ws.Cells.Item(354, 3).Value = $ws.Cells.Item(353, 3).Value2
$ws.Cells.Item(354, 5).Value = $ws.Cells.Item(353, 5).Value2
$ws.Cells.Item(354, 6).Value = $ws.Cells.Item(353, 6).Value2
$ws.Cells.Item(354, 7).Value = $ws.Cells.Item(353, 7).Value2
$ws.Cells.Item(354, 8).Value = $ws.Cells.Item(353, 8).Value2
$ws.Cells.Item(354, 9).Value = $ws.Cells.Item(353, 9).Value2
$ws.Cells.Item(354, 18).Value = $ws.Cells.Item(353, 18).Value2

# Row 335
$ws.Cells.Item(335, 4).Value = 44706
$ws.Cells.Item(335, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(335, 10).Value = 300
$ws.Cells.Item(335, 11).Value = 18000
$ws.Cells.Item(335, 12).Value = 18000
$ws.Cells.Item(335, 13).Value = 18000
$ws.Cells.Item(335, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(335, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(335, 16).Value = 360
$ws.Cells.Item(335, 17).Value = 50

# Row 336
$ws.Cells.Item(336, 4).Value = 44664
$ws.Cells.Item(336, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(336, 10).Value = 400
$ws.Cells.Item(336, 11).Value = 7000
$ws.Cells.Item(336, 12).Value = 7000
$ws.Cells.Item(336, 13).Value = 7000
$ws.Cells.Item(336, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(336, 15).Value = 'Región del Maule'
$ws.Cells.Item(336, 16).Value = 140
$ws.Cells.Item(336, 17).Value = 50

# Row 337
$ws.Cells.Item(337, 4).Value = 44566
$ws.Cells.Item(337, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(337, 10).Value = 300
$ws.Cells.Item(337, 11).Value = 5000
$ws.Cells.Item(337, 12).Value = 5000
$ws.Cells.Item(337, 13).Value = 5000
$ws.Cells.Item(337, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(337, 15).Value = 'Región del Maule'
$ws.Cells.Item(337, 16).Value = 83
$ws.Cells.Item(337, 17).Value = 60

# Row 338
$ws.Cells.Item(338, 4).Value = 44344
$ws.Cells.Item(338, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(338, 10).Value = 300
$ws.Cells.Item(338, 11).Value = 9000
$ws.Cells.Item(338, 12).Value = 9000
$ws.Cells.Item(338, 13).Value = 9000
$ws.Cells.Item(338, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(338, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(338, 16).Value = 180
$ws.Cells.Item(338, 17).Value = 50

# Row 339
$ws.Cells.Item(339, 4).Value = 44351
$ws.Cells.Item(339, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(339, 10).Value = 400
$ws.Cells.Item(339, 11).Value = 8000
$ws.Cells.Item(339, 12).Value = 8000
$ws.Cells.Item(339, 13).Value = 8000
$ws.Cells.Item(339, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(339, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(339, 16).Value = 160
$ws.Cells.Item(339, 17).Value = 50

# Row 340
$ws.Cells.Item(340, 4).Value = 44508
$ws.Cells.Item(340, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(340, 10).Value = 400
$ws.Cells.Item(340, 11).Value = 8000
$ws.Cells.Item(340, 12).Value = 8000
$ws.Cells.Item(340, 13).Value = 8000
$ws.Cells.Item(340, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(340, 15).Value = 'Región del Maule'
$ws.Cells.Item(340, 16).Value = 133
$ws.Cells.Item(340, 17).Value = 60

# Row 341
$ws.Cells.Item(341, 4).Value = 44600
$ws.Cells.Item(341, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(341, 10).Value = 400
$ws.Cells.Item(341, 11).Value = 6000
$ws.Cells.Item(341, 12).Value = 6000
$ws.Cells.Item(341, 13).Value = 6000
$ws.Cells.Item(341, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(341, 15).Value = 'Región del Maule'
$ws.Cells.Item(341, 16).Value = 120
$ws.Cells.Item(341, 17).Value = 50

# Row 342
$ws.Cells.Item(342, 4).Value = 44323
$ws.Cells.Item(342, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(342, 10).Value = 300
$ws.Cells.Item(342, 11).Value = 9000
$ws.Cells.Item(342, 12).Value = 9000
$ws.Cells.Item(342, 13).Value = 9000
$ws.Cells.Item(342, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(342, 15).Value = 'Región del Maule'
$ws.Cells.Item(342, 16).Value = 150
$ws.Cells.Item(342, 17).Value = 60

# Row 343
$ws.Cells.Item(343, 4).Value = 44515
$ws.Cells.Item(343, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(343, 10).Value = 400
$ws.Cells.Item(343, 11).Value = 6000
$ws.Cells.Item(343, 12).Value = 6000
$ws.Cells.Item(343, 13).Value = 6000
$ws.Cells.Item(343, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(343, 15).Value = 'Región del Maule'
$ws.Cells.Item(343, 16).Value = 100
$ws.Cells.Item(343, 17).Value = 60

# Row 344
$ws.Cells.Item(344, 4).Value = 44602
$ws.Cells.Item(344, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(344, 10).Value = 300
$ws.Cells.Item(344, 11).Value = 7000
$ws.Cells.Item(344, 12).Value = 7000
$ws.Cells.Item(344, 13).Value = 7000
$ws.Cells.Item(344, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(344, 15).Value = 'Región del Maule'
$ws.Cells.Item(344, 16).Value = 140
$ws.Cells.Item(344, 17).Value = 50

# Row 345
$ws.Cells.Item(345, 4).Value = 44326
$ws.Cells.Item(345, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(345, 10).Value = 300
$ws.Cells.Item(345, 11).Value = 9000
$ws.Cells.Item(345, 12).Value = 9000
$ws.Cells.Item(345, 13).Value = 9000
$ws.Cells.Item(345, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(345, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(345, 16).Value = 180
$ws.Cells.Item(345, 17).Value = 50

# Row 346
$ws.Cells.Item(346, 4).Value = 44326
$ws.Cells.Item(346, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(346, 10).Value = 400
$ws.Cells.Item(346, 11).Value = 8000
$ws.Cells.Item(346, 12).Value = 8000
$ws.Cells.Item(346, 13).Value = 8000
$ws.Cells.Item(346, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(346, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(346, 16).Value = 133
$ws.Cells.Item(346, 17).Value = 60

# Row 347
$ws.Cells.Item(347, 4).Value = 44165
$ws.Cells.Item(347, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(347, 10).Value = 500
$ws.Cells.Item(347, 11).Value = 4500
$ws.Cells.Item(347, 12).Value = 4500
$ws.Cells.Item(347, 13).Value = 4500
$ws.Cells.Item(347, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(347, 15).Value = 'Región del Maule'
$ws.Cells.Item(347, 16).Value = 75
$ws.Cells.Item(347, 17).Value = 60

# Row 348
$ws.Cells.Item(348, 4).Value = 44655
$ws.Cells.Item(348, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(348, 10).Value = 500
$ws.Cells.Item(348, 11).Value = 6000
$ws.Cells.Item(348, 12).Value = 6000
$ws.Cells.Item(348, 13).Value = 6000
$ws.Cells.Item(348, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(348, 15).Value = 'Región del Maule'
$ws.Cells.Item(348, 16).Value = 120
$ws.Cells.Item(348, 17).Value = 50

# Row 349
$ws.Cells.Item(349, 4).Value = 44315
$ws.Cells.Item(349, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(349, 10).Value = 300
$ws.Cells.Item(349, 11).Value = 7000
$ws.Cells.Item(349, 12).Value = 7000
$ws.Cells.Item(349, 13).Value = 7000
$ws.Cells.Item(349, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(349, 15).Value = 'Región del Maule'
$ws.Cells.Item(349, 16).Value = 117
$ws.Cells.Item(349, 17).Value = 60

# Row 350
$ws.Cells.Item(350, 4).Value = 44448
$ws.Cells.Item(350, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(350, 10).Value = 300
$ws.Cells.Item(350, 11).Value = 16000
$ws.Cells.Item(350, 12).Value = 16000
$ws.Cells.Item(350, 13).Value = 16000
$ws.Cells.Item(350, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(350, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(350, 16).Value = 320
$ws.Cells.Item(350, 17).Value = 50

# Row 351
$ws.Cells.Item(351, 4).Value = 44263
$ws.Cells.Item(351, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(351, 10).Value = 300
$ws.Cells.Item(351, 11).Value = 5000
$ws.Cells.Item(351, 12).Value = 5000
$ws.Cells.Item(351, 13).Value = 5000
$ws.Cells.Item(351, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(351, 15).Value = 'Región del Maule'
$ws.Cells.Item(351, 16).Value = 83
$ws.Cells.Item(351, 17).Value = 60

# Row 352
$ws.Cells.Item(352, 4).Value = 44648
$ws.Cells.Item(352, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(352, 10).Value = 500
$ws.Cells.Item(352, 11).Value = 7000
$ws.Cells.Item(352, 12).Value = 7000
$ws.Cells.Item(352, 13).Value = 7000
$ws.Cells.Item(352, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(352, 15).Value = 'Región del Maule'
$ws.Cells.Item(352, 16).Value = 140
$ws.Cells.Item(352, 17).Value = 50

# Row 353
$ws.Cells.Item(353, 4).Value = 44376
$ws.Cells.Item(353, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(353, 10).Value = 400
$ws.Cells.Item(353, 11).Value = 7500
$ws.Cells.Item(353, 12).Value = 7500
$ws.Cells.Item(353, 13).Value = 7500
$ws.Cells.Item(353, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(353, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(353, 16).Value = 150
$ws.Cells.Item(353, 17).Value = 50

# Row 354
$ws.Cells.Item(354, 4).Value = 44442
$ws.Cells.Item(354, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(354, 10).Value = 250
$ws.Cells.Item(354, 11).Value = 15000
$ws.Cells.Item(354, 12).Value = 15000
$ws.Cells.Item(354, 13).Value = 15000
$ws.Cells.Item(354, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(354, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(354, 16).Value = 300
$ws.Cells.Item(354, 17).Value = 50
